$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.761.71"
$ws.Cells.Item(2, 5).Value = "  -2.84%  "
$ws.Cells.Item(3, 4).Value = "1.783.33"
$ws.Cells.Item(3, 5).Value = "  -3.03%  "
$ws.Cells.Item(4, 4).Value = "1.004"
$ws.Cells.Item(4, 5).Value = "  +0.24%  "
$ws.Cells.Item(5, 4).Value = "242.11"
$ws.Cells.Item(5, 5).Value = "  -7.24%  "
$ws.Cells.Item(6, 4).Value = "1.004"
$ws.Cells.Item(6, 5).Value = "  +0.32%  "
$ws.Cells.Item(7, 4).Value = "0.5072"
$ws.Cells.Item(7, 5).Value = "  -3.45%  "
$ws.Cells.Item(8, 4).Value = "42.39"
$ws.Cells.Item(8, 5).Value = "  -5.28%  "
$ws.Cells.Item(9, 4).Value = "0.2516"
$ws.Cells.Item(9, 5).Value = "  -21.21%  "
$ws.Cells.Item(10, 4).Value = "0.06195"
$ws.Cells.Item(10, 5).Value = "  -8.81%  "
$ws.Cells.Item(11, 4).Value = "1.811.70"
$ws.Cells.Item(11, 5).Value = "  -1.52%  "
$ws.Cells.Item(12, 4).Value = "0.06855"
$ws.Cells.Item(12, 5).Value = "  -11.61%  "
$ws.Cells.Item(13, 5).Value = "  -20.32%  "
$ws.Cells.Item(14, 4).Value = "0.6200"
$ws.Cells.Item(14, 5).Value = "  -20.93%  "
$ws.Cells.Item(15, 4).Value = "78.63"
$ws.Cells.Item(15, 5).Value = "  -10.52%  "
$ws.Cells.Item(16, 4).Value = "4.426"
$ws.Cells.Item(16, 5).Value = "  -11.71%  "
$ws.Cells.Item(17, 5).Value = "  +0.30%  "
$ws.Cells.Item(18, 4).Value = "1.004"
$ws.Cells.Item(18, 5).Value = "  +0.32%  "
$ws.Cells.Item(19, 4).Value = "25.764.51"
$ws.Cells.Item(19, 5).Value = "  -2.95%  "
$ws.Cells.Item(20, 4).Value = "11.52"
$ws.Cells.Item(20, 5).Value = "  -16.79%  "
$ws.Cells.Item(21, 4).Value = "2.038.42"
$ws.Cells.Item(21, 5).Value = "  -1.53%  "
$ws.Cells.Item(22, 4).Value = "0.000006338"
$ws.Cells.Item(22, 5).Value = "  -20.27%  "
$ws.Cells.Item(23, 4).Value = "3.963"
$ws.Cells.Item(23, 5).Value = "  -14.29%  "
$ws.Cells.Item(24, 4).Value = "5.239"
$ws.Cells.Item(24, 5).Value = "  -12.29%  "
$ws.Cells.Item(25, 4).Value = "8.112"
$ws.Cells.Item(25, 5).Value = "  -13.13%  "
$ws.Cells.Item(26, 4).Value = "131.82"
$ws.Cells.Item(26, 5).Value = "  -7.08%  "
$ws.Cells.Item(27, 4).Value = "1.911"
$ws.Cells.Item(27, 5).Value = "  -13.00%  "
$ws.Cells.Item(28, 4).Value = "14.58"
$ws.Cells.Item(28, 5).Value = "  -13.97%  "
$ws.Cells.Item(29, 4).Value = "1.382"
$ws.Cells.Item(29, 5).Value = "  -17.67%  "
$ws.Cells.Item(30, 4).Value = "99.44"
$ws.Cells.Item(30, 5).Value = "  -10.88%  "
$ws.Cells.Item(31, 4).Value = "0.08331"
$ws.Cells.Item(31, 5).Value = "  -4.29%  "
$ws.Cells.Item(32, 4).Value = "3.607"
$ws.Cells.Item(33, 4).Value = "0.04343"
$ws.Cells.Item(33, 5).Value = "  -10.99%  "
$ws.Cells.Item(34, 4).Value = "2.737"
$ws.Cells.Item(34, 5).Value = "  -4.33%  "
$ws.Cells.Item(35, 4).Value = "3.174"
$ws.Cells.Item(35, 5).Value = "  -22.13%  "
$ws.Cells.Item(36, 4).Value = "1.044"
$ws.Cells.Item(36, 5).Value = "  -8.05%  "
$ws.Cells.Item(37, 4).Value = "0.6287"
$ws.Cells.Item(37, 5).Value = "  -13.59%  "
$ws.Cells.Item(38, 4).Value = "2.828"
$ws.Cells.Item(38, 5).Value = "  -8.59%  "
$ws.Cells.Item(39, 4).Value = "2.098"
$ws.Cells.Item(39, 5).Value = "  -6.30%  "
$ws.Cells.Item(40, 4).Value = "1.005"
$ws.Cells.Item(40, 5).Value = "  +0.37%  "
$ws.Cells.Item(41, 2).Value = "Quant"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(41, 4).Value = "100.89"
$ws.Cells.Item(41, 5).Value = "  -7.88%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42, 4).Value = "0.01460"
$ws.Cells.Item(42, 5).Value = "  -16.80%  "
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).Value = "0.7817"
$ws.Cells.Item(43, 5).Value = "  -12.71%  "
$ws.Cells.Item(44, 4).Value = "0.3929"
$ws.Cells.Item(44, 5).Value = "  -18.24%  "
$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(45, 4).Value = "5.202"
$ws.Cells.Item(45, 5).Value = "  -12.38%  "
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).Value = "6.240"
$ws.Cells.Item(46, 5).Value = "  -18.61%  "
$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(47, 4).Value = "0.05257"
$ws.Cells.Item(47, 5).Value = "  -10.06%  "
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).Value = "52.81"
$ws.Cells.Item(48, 5).Value = "  -11.47%  "
$ws.Cells.Item(49, 2).Value = "USDD"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(49, 4).Value = "1.007"
$ws.Cells.Item(49, 5).Value = "  +0.65%  "
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "29.39"
$ws.Cells.Item(50, 5).Value = "  -15.83%  "
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "7.499"
$ws.Cells.Item(51, 5).Value = "  -16.33%  "
